# Update cryptos list price (D) and volume change (E) values per upstream data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.238.30"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "1.834.70"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9984"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.17"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6198"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07377"
$ws.Range("E8").Value = "  -1.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2923"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.34"
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07674"
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("D12").Value = "1.825.25"
$ws.Range("E12").Value = "  -0.83%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.982"
$ws.Range("E13").Value = "  -0.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6714"
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.75"
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008963"
$ws.Range("E16").Value = "  -3.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.874"
$ws.Range("E17").Value = "  -1.14%  "
$ws.Range("D18").Value = "29.211.78"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").Value = "2.070.24"
$ws.Range("E19").Value = "  -1.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "236.43"
$ws.Range("E20").Value = "  +2.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.52"
$ws.Range("E21").Value = "  -1.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.360"
$ws.Range("E23").Value = "  +2.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9995"
$ws.Range("E24").Value = "  -0.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.31"
$ws.Range("E25").Value = "  -1.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1402"
$ws.Range("E26").Value = "  +0.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.570"
$ws.Range("E27").Value = "  +0.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.66"
$ws.Range("E28").Value = "  -1.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.493"
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05766"
$ws.Range("E30").Value = "  +4.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.114"
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.092"
$ws.Range("E32").Value = "  -1.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.209"
$ws.Range("E33").Value = "  +0.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.855"
$ws.Range("E34").Value = "  +0.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7330"
$ws.Range("E35").Value = "  -1.66%  "
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.608"
$ws.Range("E37").Value = "  -2.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.861"
$ws.Range("E38").Value = "  +3.34%  "
$ws.Range("D39").Value = "1.226.68"
$ws.Range("E39").Value = "  +0.55%  "
$ws.Range("E40").Value = "  -1.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.262"
$ws.Range("E41").Value = "  -3.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9120"
$ws.Range("E42").Value = "  +1.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.75"
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("D45").Value = "1.973.88"
$ws.Range("E45").Value = "  -1.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.61"
$ws.Range("E46").Value = "  -0.32%  "
$ws.Range("E47").Value = "  -1.31%  "
$ws.Range("E48").Value = "  -2.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.175"
$ws.Range("E49").Value = "  +0.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4031"
$ws.Range("E50").Value = "  -1.03%  "
